$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.422.25"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.797.57"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'223.96"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").Value = "'0.599"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'39.23"
$ws.Range("E8").Value = "  +7.15%  "
$ws.Range("D9").Value = "'0.286"
$ws.Range("E9").Value = "  -4.97%  "
$ws.Range("D10").Value = "'0.0665"
$ws.Range("E10").Value = "  -5.10%  "
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("D12").Value = "2.058.05"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "'10.82"
$ws.Range("E13").Value = "  -5.99%  "
$ws.Range("D14").Value = "1.794.72"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "34.395.33"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "'0.625"
$ws.Range("E16").Value = "  -4.45%  "
$ws.Range("D17").Value = "'4.34"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").Value = "'67.83"
$ws.Range("E18").Value = "  -3.13%  "
$ws.Range("D19").Value = "'238.51"
$ws.Range("E19").Value = "  -3.10%  "
$ws.Range("D20").Value = "0.0₃0760"
$ws.Range("E20").Value = "  -4.40%  "
$ws.Range("D21").Value = "'11.05"
$ws.Range("E21").Value = "  -4.79%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "'4.06"
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("D24").Value = "'2.15"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").Value = "'170.47"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").Value = "'17.54"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").Value = "'7.63"
$ws.Range("E27").Value = "  -4.13%  "
$ws.Range("D28").Value = "'0.120"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("D31").Value = "'3.72"
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("E32").Value = "  -3.59%  "
$ws.Range("D33").Value = "'3.82"
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("E36").Value = "  -5.19%  "
$ws.Range("D37").Value = "1.297.39"
$ws.Range("E37").Value = "  -6.99%  "
$ws.Range("D38").Value = "'0.0184"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").Value = "'2.29"
$ws.Range("E39").Value = "  -6.71%  "
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("D42").Value = "'81.30"
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("D44").Value = "'0.936"
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("D45").Value = "'13.98"
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("E46").Value = "  +4.44%  "
$ws.Range("D47").Value = "1.959.38"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  -5.76%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").Value = "'101.22"
$ws.Range("E50").Value = "  -3.00%  "
$ws.Range("E51").Value = "  -0.61%  "
